# 16th commit: refresh the demo web-login test credentials and leave the
# DemoWebLogin sheet as the active tab (matches the author's workflow of
# pasting a freshly generated email/password pair before saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DemoWebLogin")

$ws.Range("A2").Value = "abcdefghijkluyiooip@gmail.com"
$ws.Range("B2").Value = "01bXbbccb"

$ws.Activate()
$ws.Range("B2").Select()
